$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: add new issue entry (TESLA/ECDSA) ---
$ws.Range("B15").Value = "Build a conditional statement to set TESLA parameters to null if ECDSA is being used"
$ws.Range("C15").Value = "ConfigParameters"

# Apply the "Good" cell style to the existing (used) cells of rows 3, 4 and 15,
# matching the green/yellow banding already used elsewhere in the sheet.
$ws.Range("A3:C3").Style = "Neutral"
$ws.Range("A4:C4").Style = "Good"
$ws.Range("A15:C15").Style = "Good"

# --- Rows 16-100: continue the Issue # numbering sequence in column A ---
for ($r = 16; $r -le 100; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# --- Update the active selection to reflect where editing left off ---
$ws.Range("B18").Select()
